$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Rows 16-22 hold the "Estado de Cuenta" worker listing (Tipo Doc, N Doc, Nombre,
# Periodo Mora, Salario Basico, Valor Mora). The workers were reshuffled into new
# row positions and one worker's overdue value (Valor Mora) was updated.
#
# Columns: C = N Doc Trabajador, D = Nombre Trabajador, F = Salario Basico, G = Valor Mora

$ws.Range("C16").Value = "30894870"
$ws.Range("D16").Value = "KAREN CECILIA SUAREZ ESALAS"
$ws.Range("F16").Value = 48000
$ws.Range("G16").Value = 1200000

$ws.Range("C17").Value = "73201827"
$ws.Range("D17").Value = "LENDER RAFAEL OROZCO BLANCO"
$ws.Range("F17").Value = 120000
$ws.Range("G17").Value = 2652000

$ws.Range("C18").Value = "73560182"
$ws.Range("D18").Value = "DAIRO ENRIQUE ROMERO HERRERA"
$ws.Range("F18").Value = 180000
$ws.Range("G18").Value = 4500000

$ws.Range("C19").Value = "64727108"
$ws.Range("D19").Value = "BRIGIDA PATRICIA VILLADIEGO MENDOZA"
$ws.Range("F19").Value = 72000
$ws.Range("G19").Value = 1800000

$ws.Range("C20").Value = "45591343"
$ws.Range("D20").Value = "YINCETH GONZALEZ GONZALEZ"
$ws.Range("F20").Value = 80000
$ws.Range("G20").Value = 908526

$ws.Range("C21").Value = "9149147"
$ws.Range("D21").Value = "ANTONIO TADEO DURAN PATERNINA"
$ws.Range("F21").Value = 100000
$ws.Range("G21").Value = 2500000

$ws.Range("C22").Value = "1143410300"
$ws.Range("D22").Value = "NICOLAS FELIPE CABARCAS CUEVAS"
$ws.Range("F22").Value = 12000
$ws.Range("G22").Value = 1500000
